$ws = $excel.ActiveWorkbook.ActiveSheet

$ws.Range('D2').Value = '62.525.48'
$ws.Range('E2').Value = '  -0.39%  '
$ws.Range('D3').Value = '2.439.56'
$ws.Range('E3').Value = '  -0.66%  '
$ws.Range('E4').Value = '  -0.17%  '
$ws.Range('E5').Value = '  +0.73%  '
$ws.Range('E6').Value = '  -1.62%  '
$ws.Range('E7').Value = '  -0.01%  '
$ws.Range('E8').Value = '  -0.58%  '
$ws.Range('D9').Value = '2.436.37'
$ws.Range('E9').Value = '  -0.78%  '
$ws.Range('E10').Value = '  -3.18%  '
$ws.Range('E11').Value = '  +0.74%  '
$ws.Range('E12').Value = '  -0.26%  '
$ws.Range('E13').Value = '  -1.23%  '
$ws.Range('E14').Value = '  -0.97%  '
$ws.Range('E15').Value = '  -2.09%  '
$ws.Range('D16').Value = '2.875.91'
$ws.Range('E16').Value = '  -1.10%  '
$ws.Range('D17').Value = '62.292.40'
$ws.Range('E17').Value = '  -0.87%  '
$ws.Range('D18').Value = '2.418.06'
$ws.Range('E18').Value = '  -1.33%  '
$ws.Range('E19').Value = '  -1.93%  '
$ws.Range('E20').Value = '  -1.37%  '
$ws.Range('E21').Value = '  +1.49%  '
$ws.Range('E22').Value = '  -0.31%  '
$ws.Range('E23').Value = '  +4.16%  '
$ws.Range('E24').Value = '  +0.33%  '
$ws.Range('E25').Value = '  -1.05%  '
$ws.Range('E26').Value = '  +2.22%  '
$ws.Range('E27').Value = '  +6.92%  '
$ws.Range('D28').Value = '0.0₃0963'
$ws.Range('E28').Value = '  -5.30%  '
$ws.Range('E30').Value = '  +0.09%  '
$ws.Range('E31').Value = '  -2.88%  '
$ws.Range('E32').Value = '  -2.21%  '
$ws.Range('E33').Value = '  -0.69%  '
$ws.Range('E34').Value = '  -4.01%  '
$ws.Range('E35').Value = '  -1.02%  '
$ws.Range('E36').Value = '  +0.18%  '
$ws.Range('E37').Value = '  -1.45%  '
$ws.Range('E38').Value = '  -1.65%  '
$ws.Range('E39').Value = '  -0.80%  '
$ws.Range('E40').Value = '  -2.91%  '
$ws.Range('E41').Value = '  +1.08%  '
$ws.Range('E42').Value = '  -2.31%  '
$ws.Range('E43').Value = '  +1.59%  '
$ws.Range('E44').Value = '  +0.01%  '
$ws.Range('E45').Value = '  -4.51%  '
$ws.Range('E46').Value = '  -1.25%  '
$ws.Range('E47').Value = '  -0.20%  '
$ws.Range('E48').Value = '  -2.24%  '
$ws.Range('E49').Value = '  -0.44%  '
$ws.Range('E50').Value = '  -4.74%  '
$ws.Range('E51').Value = '  -1.77%  '

# Numeric-looking price strings must stay TEXT (matching the original
# inlineStr cells) instead of being auto-converted to numbers by Excel.
# Route them through a text formula, then Copy/PasteSpecial-values to
# flatten to a static value without touching NumberFormat/style.
$ws.Range('D5').Formula = '="573.99"'
$ws.Range('D5').Copy() | Out-Null
$ws.Range('D5').PasteSpecial(-4163) | Out-Null
$ws.Range('D6').Formula = '="144.25"'
$ws.Range('D6').Copy() | Out-Null
$ws.Range('D6').PasteSpecial(-4163) | Out-Null
$ws.Range('D12').Formula = '="5.23"'
$ws.Range('D12').Copy() | Out-Null
$ws.Range('D12').PasteSpecial(-4163) | Out-Null
$ws.Range('D14').Formula = '="26.58"'
$ws.Range('D14').Copy() | Out-Null
$ws.Range('D14').PasteSpecial(-4163) | Out-Null
$ws.Range('D19').Formula = '="11.12"'
$ws.Range('D19').Copy() | Out-Null
$ws.Range('D19').PasteSpecial(-4163) | Out-Null
$ws.Range('D20').Formula = '="7.17"'
$ws.Range('D20').Copy() | Out-Null
$ws.Range('D20').PasteSpecial(-4163) | Out-Null
$ws.Range('D21').Formula = '="327.82"'
$ws.Range('D21').Copy() | Out-Null
$ws.Range('D21').PasteSpecial(-4163) | Out-Null
$ws.Range('D25').Formula = '="65.46"'
$ws.Range('D25').Copy() | Out-Null
$ws.Range('D25').PasteSpecial(-4163) | Out-Null
$ws.Range('D26').Formula = '="633.21"'
$ws.Range('D26').Copy() | Out-Null
$ws.Range('D26').PasteSpecial(-4163) | Out-Null
$ws.Range('D27').Formula = '="9.18"'
$ws.Range('D27').Copy() | Out-Null
$ws.Range('D27').PasteSpecial(-4163) | Out-Null
$ws.Range('D32').Formula = '="8.07"'
$ws.Range('D32').Copy() | Out-Null
$ws.Range('D32').PasteSpecial(-4163) | Out-Null
$ws.Range('D33').Formula = '="1.88"'
$ws.Range('D33').Copy() | Out-Null
$ws.Range('D33').PasteSpecial(-4163) | Out-Null
$ws.Range('D34').Formula = '="0.136"'
$ws.Range('D34').Copy() | Out-Null
$ws.Range('D34').PasteSpecial(-4163) | Out-Null
$ws.Range('D40').Formula = '="5.24"'
$ws.Range('D40').Copy() | Out-Null
$ws.Range('D40').PasteSpecial(-4163) | Out-Null
$ws.Range('D41').Formula = '="146.17"'
$ws.Range('D41').Copy() | Out-Null
$ws.Range('D41').PasteSpecial(-4163) | Out-Null
$ws.Range('D46').Formula = '="145.41"'
$ws.Range('D46').Copy() | Out-Null
$ws.Range('D46').PasteSpecial(-4163) | Out-Null
$ws.Range('D48').Formula = '="0.0524"'
$ws.Range('D48').Copy() | Out-Null
$ws.Range('D48').PasteSpecial(-4163) | Out-Null
$ws.Range('D49').Formula = '="0.598"'
$ws.Range('D49').Copy() | Out-Null
$ws.Range('D49').PasteSpecial(-4163) | Out-Null
$ws.Range('D50').Formula = '="19.68"'
$ws.Range('D50').Copy() | Out-Null
$ws.Range('D50').PasteSpecial(-4163) | Out-Null

$excel.CutCopyMode = 0
